# Auto commit at 2025-12-09  8:05:40.75
# Updates the "Metrics" sheet values (B2:B13), which cascade via formulas
# into the "today" sheet (B11:B22, E11:E22, F11:F22, and A1's TODAY()-1),
# and restores each sheet's saved cursor/selection cell.

$wb = $excel.ActiveWorkbook

# --- Metrics sheet: refreshed metric values ---
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value  = 108159.6
$metrics.Range("B3").Value  = 92985.010000000009
$metrics.Range("B4").Value  = 33244.089999999997
$metrics.Range("B5").Value  = 4381
$metrics.Range("B6").Value  = 5310866.7100000009
$metrics.Range("B7").Value  = 4493337.9700000007
$metrics.Range("B8").Value  = 1565200.9700000004
$metrics.Range("B9").Value  = 207088
$metrics.Range("B10").Value = 33776247.699999996
$metrics.Range("B11").Value = 31768613.130000003
$metrics.Range("B12").Value = 11846923.009999996
$metrics.Range("B13").Value = 1304718

# --- "today" sheet: formulas pull from Metrics automatically; A1's
# TODAY()-1 recalculates from the runtime clock. Just restore selection. ---
$today = $wb.Worksheets.Item("today")

# Move the cursor on "today" first, then back to "Metrics" last so the
# final ActiveSheet/selection state matches the saved workbook (the
# "today" tab is the one marked tabSelected="1").
$today.Range("G7").Select()
$metrics.Range("F21").Select()
$today.Range("G7").Select()
